# FindNexialCommands.xlsx - add "storeKeys(json,jsonpath,var)" to the json
# command list and drop the stray "text" stub column from the "#system"
# lookup sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1) json command list (column M / #13): insert "storeKeys(...)" before
#        "storeValue(...)" at row 16, pushing the remaining two rows down.
$ws.Cells.Item(18, 13).Value = $ws.Cells.Item(17, 13).Value2
$ws.Cells.Item(17, 13).Value = $ws.Cells.Item(16, 13).Value2
$ws.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"

# --- 2) target list (column A / #1): remove the "text" entry (row 25) and
#        pull the remaining rows (web..xml) up by one.
for ($r = 25; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r + 1, 1).Value2
}
$ws.Cells.Item(31, 1).ClearContents()

# --- 3) drop the unused single-column "text" stub (column Y / #25): shift
#        columns Z:AE (web, webalert, webcookie, ws, ws.async, xml) one
#        column to the left, into Y:AD, for every row.
for ($r = 1; $r -le 129; $r++) {
    for ($c = 26; $c -le 31; $c++) {
        $ws.Cells.Item($r, $c - 1).Value = $ws.Cells.Item($r, $c).Value2
    }
    $ws.Cells.Item($r, 31).ClearContents()
}

# --- 4) fix up the defined names so they still point at the right ranges.
$wb.Names.Item("json").RefersTo      = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$27"
